$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 36: date (06.03.2018 -> serial 43165), task description, duration
# Copy the formatting from the row above (A35) so the new date cell picks up
# the same style (numFmtId 14, "date" style) instead of minting a new xf.
$ws.Range("A35").Copy() | Out-Null
$ws.Range("A36").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A36").Value = 43165

$ws.Range("B35").Copy() | Out-Null
$ws.Range("B36").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B36").Value = "Création de la page qui permettra d'ajouter des nouveaux articles "

$ws.Range("C36").Value = "1 période"

$excel.CutCopyMode = 0

# Update selection to match the post-edit state (B37 selected)
$ws.Range("B37").Select()
